# Updates cryptos list (Price/Volume(1h) columns, plus a row 35/36 swap
# for Filecoin <-> PEPE) per the "Updated cryptos list" GitHub Actions
# commit. Numeric-looking Price strings are written with a leading
# apostrophe (Excel's quote-prefix convention) so they stay text cells
# instead of being auto-coerced to numbers, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.222.98"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "2.978.98"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D5").Value = "'566.00"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "'136.94"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "2.974.25"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  +7.82%  "
$ws.Range("D12").Value = "'0.449"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'33.54"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "3.468.20"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "'7.04"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "2.971.74"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "59.191.41"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").Value = "'435.11"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").Value = "'13.63"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'0.723"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "'7.01"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "'13.05"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("D25").Value = "'79.85"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'2.22"
$ws.Range("E27").Value = "  +6.24%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "'7.70"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'6.17"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("E33").Value = "  +5.45%  "
$ws.Range("D34").Value = "'0.988"
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.88"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0763"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").Value = "'2.06"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "'48.51"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").Value = "'2.78"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "'395.01"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "2.715.26"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D47").Value = "'122.37"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").Value = "'34.41"
$ws.Range("E48").Value = "  +11.85%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  -2.35%  "
$ws.Range("D51").Value = "'23.17"
$ws.Range("E51").Value = "  -0.97%  "
